$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.032.22"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.79"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.57"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.11"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.566.70"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.520"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.053.81"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.87"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.97"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.94"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.05"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E30").Value = "  +3.99%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  +4.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.427.05"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("E34").Value = "  +12.52%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0168"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.813"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.83"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.702.92"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.67"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +0.34%  "
